$p = $ppt.ActivePresentation

# --- Slide 3 ("Magic Trick Revealed (I)"): resize the body placeholder and
#     rework the final sentence so it distinguishes the hidden/listed suits. ---
$s3 = $p.Slides.Item(3)
$rect3 = $s3.Shapes.Item("Rectangle 3")

# <a:off x="469900" y="1371600"/><a:ext cx="8216900" cy="3124200"/>
#   -> <a:off x="304800" y="1371600"/><a:ext cx="8534400" cy="3276600"/>
# Shape.Left/Top/Width/Height are in points; OOXML stores EMU (12700 EMU = 1 pt).
$rect3.Left   = 304800 / 12700
$rect3.Top    = 1371600 / 12700
$rect3.Width  = 8534400 / 12700
$rect3.Height = 3276600 / 12700

# Paragraph 3 of that placeholder reads "A lists one of them 1st".
# Run 2 is " lists one of them " - replace it with the longer phrasing while
# keeping a leading space-only run (same rPr) ahead of the new text, exactly
# like the authored edit split it into two runs.
$tr3 = $rect3.TextFrame.TextRange
$para3 = $tr3.Paragraphs(3)
$run2 = $para3.Runs(2)
$run2.Text = " hides one, lists the other one "

# The host only materialises a run split once a character-level format is
# touched; re-apply the (unchanged) font size of the first character so the
# leading " " becomes its own run, matching the other run's formatting.
$para3b = $tr3.Paragraphs(3)
$run2b = $para3b.Runs(2)
$firstChar = $tr3.Characters($run2b.Start, 1)
$firstChar.Font.Size = $firstChar.Font.Size
